# Add two new "factor_aesthetics" groups (marine_system, climate_threat)
# below the existing data, matching columns:
#   A = variable, B = level, C = label, D = order, E = colour
#
# The cell-write order below is deliberately chosen (rather than a simple
# left-to-right / row-by-row fill) so that new values land in the shared
# string table in the same sequence as the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "marine_system"
$ws.Range("C16").Value = "Coastal land"
$ws.Range("C17").Value = "Coastal ocean"
$ws.Range("C18").Value = "Open-ocean"
$ws.Range("A17").Value = "marine_system"
$ws.Range("A18").Value = "marine_system"
$ws.Range("A19").Value = "climate_threat"
$ws.Range("B19").Value = "Temperature"
$ws.Range("C19").Value = "Temperature"
$ws.Range("B20").Value = "SLR"
$ws.Range("B21").Value = "Extreme_weather"
$ws.Range("E16").Value = "NA"
$ws.Range("E17").Value = "NA"
$ws.Range("E18").Value = "NA"
$ws.Range("E19").Value = "NA"
$ws.Range("E20").Value = "NA"
$ws.Range("E21").Value = "NA"
$ws.Range("E22").Value = "NA"
$ws.Range("C20").Value = "Sea level rise"
$ws.Range("C21").Value = "Extreme weather"
$ws.Range("B16").Value = "land"
$ws.Range("B18").Value = "open_ocean"
$ws.Range("B17").Value = "coastal_ocean"
$ws.Range("B22").Value = "NA"
$ws.Range("C22").Value = "Unidentified"
$ws.Range("A20").Value = "climate_threat"
$ws.Range("A21").Value = "climate_threat"
$ws.Range("A22").Value = "climate_threat"
$ws.Range("D16").Value = 1
$ws.Range("D17").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("D19").Value = 1
$ws.Range("D20").Value = 2
$ws.Range("D21").Value = 3
$ws.Range("D22").Value = 4

# Scroll the view down a bit and leave the selection on C22, matching
# where editing left off.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("C22").Select()
